# CajaInternet.xlsx — "segundo dia" update
# Adds the second day's cash-drawer count (opening date 10-Jun-2018 in the
# second table) and fills in the "Finalice" / closing counts for all three
# coin-count tables on the first block (rows 4-13), recalculating totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Second table's date header (G2) — second day, one week after A2 (03-Jun-2018)
$ws.Range("G2").Value = 43261

# Row 4  (0.05)
$ws.Range("D4").Value = 0
$ws.Range("H4").Value = 3
$ws.Range("J4").Value = 3

# Row 5  (0.10)
$ws.Range("D5").Value = 1
$ws.Range("H5").Value = 18
$ws.Range("J5").Value = 18

# Row 6  (0.25)
$ws.Range("D6").Value = 1
$ws.Range("H6").Value = 7
$ws.Range("J6").Value = 5

# Row 7  (0.50)
$ws.Range("D7").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("J7").Value = 0

# Row 8  (1)
$ws.Range("D8").Value = 9
$ws.Range("H8").Value = 6
$ws.Range("J8").Value = 0

# Row 9  (5)
$ws.Range("D9").Value = 0
$ws.Range("H9").Value = 10
$ws.Range("J9").Value = 11

# Row 10 (10)
$ws.Range("D10").Value = 10
$ws.Range("H10").Value = 2
$ws.Range("J10").Value = 6

# Row 11 (20)
$ws.Range("D11").Value = 3
$ws.Range("H11").Value = 1
$ws.Range("J11").Value = 3

# Row 12 (50)
$ws.Range("D12").Value = 2
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0

# Row 13 (100)
$ws.Range("D13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0

# Grand total for the "Finalice" column of the first table (E14) — the other
# two totals (I14/K14) already have this SUM formula and just pick up the
# new inputs above.
$ws.Range("E14").Formula = "=SUM(E4:E13)"

# Leave the selection where the author last left it when saving.
$ws.Range("J14").Select()
